$d = $word.ActiveDocument

function Insert-SectionPrefix($searchText, $prefixText) {
    $rng = $d.Content
    $rng.Find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

    $prefixRange = $rng.Duplicate
    $prefixRange.Collapse(1)
    $prefixRange.Text = $prefixText

    # Re-grab the just-inserted range and round-trip its FormattedText through
    # itself; this forces the run boundary between the new prefix and the
    # original heading text to be preserved as two separate runs (both with
    # identical bold/size formatting) instead of being silently coalesced
    # into a single run on save.
    $newRange = $d.Range($prefixRange.Start, $prefixRange.End)
    $newRange.FormattedText = $newRange.FormattedText
}

Insert-SectionPrefix "WORKFLOW" "SECTION ONE: "
Insert-SectionPrefix "MODIFYING FIELDS" "SECTION TWO: "
